$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.929.78'
$ws.Range("E2").Value = '  +0.22%  '

# Row 3
$ws.Range("D3").Value = '3.161.64'
$ws.Range("E3").Value = '  +0.10%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.52'
$ws.Range("E5").Value = '  +1.10%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.31'
$ws.Range("E6").Value = '  -0.60%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").Value = '3.161.63'
$ws.Range("E8").Value = '  +0.14%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -0.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  -1.72%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.15'
$ws.Range("E11").Value = '  -0.71%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.504'
$ws.Range("E12").Value = '  -0.29%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000265'
$ws.Range("E13").Value = '  +2.41%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.44'
$ws.Range("E14").Value = '  -1.63%  '

# Row 15
$ws.Range("D15").Value = '3.676.58'
$ws.Range("E15").Value = '  +0.15%  '

# Row 16
$ws.Range("D16").Value = '64.851.73'
$ws.Range("E16").Value = '  -0.04%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.20'
$ws.Range("E17").Value = '  -0.68%  '

# Row 18
$ws.Range("D18").Value = '3.119.94'
$ws.Range("E18").Value = '  -1.33%  '

# Row 19
$ws.Range("E19").Value = '  +0.45%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '509.02'
$ws.Range("E20").Value = '  -2.03%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.00'
$ws.Range("E21").Value = '  +0.04%  '

# Row 22
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.718'
$ws.Range("E22").Value = '  -2.61%  '

# Row 23
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.31'
$ws.Range("E23").Value = '  +0.72%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.77'
$ws.Range("E24").Value = '  -1.09%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.72'
$ws.Range("E25").Value = '  -0.68%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.18%  '

# Row 27
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.97'
$ws.Range("E27").Value = '  +2.13%  '

# Row 28
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.92'
$ws.Range("E28").Value = '  -0.51%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("E29").Value = '  +0.15%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.83'
$ws.Range("E30").Value = '  +5.89%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.77'
$ws.Range("E31").Value = '  -0.96%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.997'
$ws.Range("E32").Value = '  -0.23%  '

# Row 33
$ws.Range("E33").Value = '  +0.91%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.24'
$ws.Range("E34").Value = '  +1.46%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.51'
$ws.Range("E35").Value = '  -1.24%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.77'
$ws.Range("E36").Value = '  -1.84%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0891'
$ws.Range("E37").Value = '  +2.86%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '478.42'
$ws.Range("E38").Value = '  -1.83%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0417'
$ws.Range("E39").Value = '  -1.29%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.94'
$ws.Range("E40").Value = '  -1.90%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.78'
$ws.Range("E41").Value = '  +1.25%  '

# Row 42
$ws.Range("D42").Value = '3.008.09'
$ws.Range("E42").Value = '  -3.41%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.116'
$ws.Range("E43").Value = '  -4.16%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.44'
$ws.Range("E44").Value = '  -1.43%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.283'
$ws.Range("E45").Value = '  -4.99%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.40'
$ws.Range("E46").Value = '  -2.93%  '

# Row 47
$ws.Range("D47").Value = '0.0₃0588'
$ws.Range("E47").Value = '  +1.30%  '

# Row 49
$ws.Range("E49").Value = '  -1.41%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.25'
$ws.Range("E50").Value = '  -2.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.49'
$ws.Range("E51").Value = '  +14.26%  '
